$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (L) mirroring the existing "2019" column (K):
# copy each source cell's formatting into the new cell, then set its value.
# xlPasteFormats = -4122

# Row 2 - empty bottom-border separator cell (same style as K2)
$ws.Range("K2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null

# Row 3 - header year value
$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Value = 2020

# Row 4-7 - data rows (first block)
$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").Value = 1004

$ws.Range("K5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null
$ws.Range("L5").Value = 8279

$ws.Range("K6").Copy() | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null
$ws.Range("L6").Value = 1752

$ws.Range("K7").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null
$ws.Range("L7").Value = 6527

# Row 8-9 - data rows (second block)
$ws.Range("K8").Copy() | Out-Null
$ws.Range("L8").PasteSpecial(-4122) | Out-Null
$ws.Range("L8").Value = 10324

$ws.Range("K9").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null
$ws.Range("L9").Value = 4131

# Row 10 - last data row (bottom border)
$ws.Range("K10").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null
$ws.Range("L10").Value = 6193

$excel.CutCopyMode = $false

# Update the saved cursor/selection position
$ws.Range("Q11").Select() | Out-Null
